$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 26 de Mayo de 2020 a las 23:05"

# Update Estados Unidos row (row 4)
$ws.Range("B4").Value = 1719766
$ws.Range("C4").Value = 13540
$ws.Range("D4").Value = 473350
$ws.Range("E4").Value = 1146005
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 606
$ws.Range("H4").Value = 100411

# Update Japon row (row 43)
$ws.Range("B43").Value = 16623
$ws.Range("C43").Value = 42
$ws.Range("D43").Value = 13810
$ws.Range("E43").Value = 1967
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 16
$ws.Range("H43").Value = 846

# Update Niger row (row 115)
$ws.Range("B115").Value = 952
$ws.Range("C115").Value = 1
$ws.Range("D115").Value = 796
$ws.Range("E115").Value = 93
$ws.Range("F115").Value = 0
$ws.Range("G115").Value = 1
$ws.Range("H115").Value = 63

# Togo and Cabo Verde swap positions in the ranking (Togo overtakes Cabo Verde)
# Row 142 now holds Togo's (updated) data
$ws.Range("A142").Value = "Togo"
$ws.Range("B142").Value = 391
$ws.Range("C142").Value = 5
$ws.Range("D142").Value = 177
$ws.Range("E142").Value = 201
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 13

# Row 143 now holds Cabo Verde's (previous) data
$ws.Range("A143").Value = "Cabo Verde"
$ws.Range("B143").Value = 390
$ws.Range("C143").Value = 0
$ws.Range("D143").Value = 155
$ws.Range("E143").Value = 231
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 1
$ws.Range("H143").Value = 4

# Update Guyana row (row 164)
$ws.Range("B164").Value = 139
$ws.Range("C164").Value = 2
$ws.Range("E164").Value = 66
